$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 - existing rows 12..48 shift down to 13..49
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new weekly price record
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "Femacal de La Calera"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44883
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 300000000
$ws.Range("G12").Value = "Espárragos"
$ws.Range("H12").Value = "Verde"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 1200
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1400
$ws.Range("M12").Value = 1400
$ws.Range("N12").Value = "$/kilo"
$ws.Range("O12").Value = "Provincia de Quillota"
$ws.Range("P12").Value = 1400
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"
